# This script reorders data rows 3-21 of the "Artfynd" sheet according to
# a fixed permutation (the rows themselves are unchanged in content; only
# their position in the sheet changes). Row 1 (header) and row 2 stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow  = 21
$lastCol  = 51   # column AY

# Mapping: old row number -> new row number
$map = @{
    3  = 15
    4  = 3
    5  = 4
    6  = 5
    7  = 6
    8  = 16
    9  = 7
    10 = 17
    11 = 18
    12 = 8
    13 = 9
    14 = 10
    15 = 11
    16 = 12
    17 = 13
    18 = 19
    19 = 20
    20 = 21
    21 = 14
}

# Read the whole block (rows 3-21, columns A:AY) into memory in one shot so
# that all source data is captured before anything gets overwritten.
$srcRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$srcValues = $srcRange.Value()

$rowCount = $lastRow - $firstRow + 1

# Build the destination array applying the permutation.
$destValues = New-Object 'object[,]' $rowCount, $lastCol

foreach ($oldRow in $map.Keys) {
    $newRow = $map[$oldRow]
    $srcIdx  = $oldRow - $firstRow + 1
    $destIdx = $newRow - $firstRow + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $destValues[$destIdx - 1, $c - 1] = $srcValues[$srcIdx, $c]
    }
}

# Some columns hold values that are stored as text but look like numbers
# or dates (e.g. a "Startdatum" of "2023-01-26", a count of "1", or a
# collection number like "2301261149"). If such a string is written
# through .Value as-is, Excel auto-converts it to a real date/number,
# which would change its stored type from the original text. Detect, from
# the source data, which columns contain such text-that-looks-numeric
# values and force those columns to Text format before writing, so the
# values round-trip as the original literal strings.
$numericLike = '^-?\d+(\.\d+)?$|^\d{4}-\d{2}-\d{2}$'
for ($c = 1; $c -le $lastCol; $c++) {
    $forceText = $false
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $v = $srcValues[$r - $firstRow + 1, $c]
        if ($v -is [string] -and $v -match $numericLike) {
            $forceText = $true
            break
        }
    }
    if ($forceText) {
        $ws.Range($ws.Cells.Item($firstRow, $c), $ws.Cells.Item($lastRow, $c)).NumberFormat = "@"
    }
}

# Write everything back out in a single operation.
$destRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$destRange.Value = $destValues
